# Generate Report for Handback
# Update the timestamp strings that record when the Handoff/Handback xliff
# files were (re-)generated, as plain text values (cells already carry the
# date-time number format via their existing style).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first data row.
$overview.Range("G2").Value = "2016-09-02 23:11:47"

# zh-cn sheet, row 2: Correspond Handoff Datetime / Correspond Handback DateTime
$zhcn.Range("H2").Value = "2016-09-02 23:11:42"
$zhcn.Range("K2").Value = "2016-09-02 23:12:00"

# de-de sheet, row 2: Correspond Handoff Datetime / Correspond Handback DateTime
$dede.Range("H2").Value = "2016-09-02 23:11:47"
$dede.Range("K2").Value = "2016-09-02 23:12:15"
